# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45205 (2023-10-06) to 45206 (2023-10-07).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2() -eq 45205) {
        $cell.Value = 45206
    }
}
